$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    Every cell across the three sheets that currently carries the old
#    status string needs to show the new one.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2) zh-cn sheet: populate "Latest Target File" (F) / "Latest Handback File"
#    (G) with hyperlinked file names, and set the "Latest Handback DateTime"
#    (H) now that a handback has actually happened.
# ---------------------------------------------------------------------------
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/642b11904a4942e4212ab614758b51098bea8826/e2e/a.md",
    "",
    "",
    "a.md"
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0b6bd31b464a106299597008dff0bf1b43da2649/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf",
    "",
    "",
    "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/642b11904a4942e4212ab614758b51098bea8826/e2e/a.md",
    "",
    "",
    "a.md"
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0b6bd31b464a106299597008dff0bf1b43da2649/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf",
    "",
    "",
    "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
) | Out-Null

$wsZh.Range("H2").Value = "2016-03-21 00:28:20"
$wsZh.Range("H3").Value = "2016-03-21 00:28:20"

# ---------------------------------------------------------------------------
# 3) de-de sheet: same shape of change, but the handback finished later so it
#    gets its own, later, timestamp.
# ---------------------------------------------------------------------------
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/642b11904a4942e4212ab614758b51098bea8826/e2e/a.md",
    "",
    "",
    "a.md"
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b3b13dbe6204ebdbacdbbadfc4e6c3d561091ce0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf",
    "",
    "",
    "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/642b11904a4942e4212ab614758b51098bea8826/e2e/a.md",
    "",
    "",
    "a.md"
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b3b13dbe6204ebdbacdbbadfc4e6c3d561091ce0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf",
    "",
    "",
    "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
) | Out-Null

$wsDe.Range("H2").Value = "2016-03-21 00:28:27"
$wsDe.Range("H3").Value = "2016-03-21 00:28:27"
